# Updates the cryptocurrency price/volume table (columns B-E, rows 2-51)
# to match the latest scrape used in the "Updated cryptos list" GitHub
# Actions commit. A few rows were also re-ranked, swapping their Coin
# name/link/price/volume with the neighbouring row (e.g. Dai <-> ICP).
#
# Price cells (column D) are written with a leading apostrophe so Excel
# keeps values such as "1.00"/"0.200" as literal text instead of
# coercing them to numbers and silently dropping trailing zeros; the
# cell style is then reset to Normal so the quote-prefix marker Excel
# adds for that doesn't linger on the cell (matches the source file,
# which stores these as plain text with no special style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'95.655.03"
$ws.Range("E2").Value = '  +2.43%  '
# Row 3
$ws.Range("D3").Value = "'3.628.29"
$ws.Range("E3").Value = '  +4.85%  '
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").Value = "'237.08"
$ws.Range("E5").Value = '  +0.85%  '
# Row 6
$ws.Range("D6").Value = "'659.51"
$ws.Range("E6").Value = '  +5.46%  '
# Row 7
$ws.Range("E7").Value = '  +2.25%  '
# Row 8
$ws.Range("E8").Value = '  +3.01%  '
# Row 9
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  -0.04%  '
# Row 10
$ws.Range("D10").Value = "'0.999"
$ws.Range("E10").Value = '  -0.76%  '
# Row 11
$ws.Range("D11").Value = "'3.627.76"
$ws.Range("E11").Value = '  +4.92%  '
# Row 12
$ws.Range("D12").Value = "'0.200"
$ws.Range("E12").Value = '  +0.23%  '
# Row 13
$ws.Range("D13").Value = "'42.31"
$ws.Range("E13").Value = '  -2.61%  '
# Row 14
$ws.Range("D14").Value = "'6.44"
$ws.Range("E14").Value = '  +3.35%  '
# Row 15
$ws.Range("D15").Value = "'4.303.36"
$ws.Range("E15").Value = '  +4.84%  '
# Row 16
$ws.Range("D16").Value = "'95.476.76"
$ws.Range("E16").Value = '  +2.38%  '
# Row 17
$ws.Range("E17").Value = '  +1.98%  '
# Row 18
$ws.Range("D18").Value = "'3.630.06"
$ws.Range("E18").Value = '  +4.31%  '
# Row 19
$ws.Range("E19").Value = '  -4.25%  '
# Row 20
$ws.Range("D20").Value = "'12.92"
$ws.Range("E20").Value = '  +9.17%  '
# Row 21
$ws.Range("D21").Value = "'18.04"
$ws.Range("E21").Value = '  -0.31%  '
# Row 22
$ws.Range("E22").Value = '  +4.60%  '
# Row 23
$ws.Range("D23").Value = "'0.483"
# Row 24
$ws.Range("D24").Value = "'504.99"
$ws.Range("E24").Value = '  +0.27%  '
# Row 25
$ws.Range("E25").Value = '  +7.73%  '
# Row 26
$ws.Range("D26").Value = "'6.66"
$ws.Range("E26").Value = '  -1.90%  '
# Row 27
$ws.Range("D27").Value = "'91.97"
$ws.Range("E27").Value = '  -3.08%  '
# Row 28
$ws.Range("D28").Value = "'3.823.57"
$ws.Range("E28").Value = '  +4.92%  '
# Row 29
$ws.Range("D29").Value = "'12.54"
$ws.Range("E29").Value = '  +2.77%  '
# Row 30
$ws.Range("D30").Value = "'3.10"
$ws.Range("E30").Value = '  +8.84%  '
# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'11.29"
$ws.Range("E31").Value = '  -0.36%  '
# Row 32
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = '  -0.01%  '
# Row 33
$ws.Range("E33").Value = '  -0.48%  '
# Row 34
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = '  -0.07%  '
# Row 35
$ws.Range("D35").Value = "'32.36"
$ws.Range("E35").Value = '  +9.42%  '
# Row 36
$ws.Range("D36").Value = "'0.177"
$ws.Range("E36").Value = '  -1.66%  '
# Row 37
$ws.Range("D37").Value = "'0.561"
$ws.Range("E37").Value = '  +0.68%  '
# Row 38
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").Value = "'8.12"
$ws.Range("E38").Value = '  +7.98%  '
# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = "'569.14"
$ws.Range("E39").Value = '  -0.58%  '
# Row 40
$ws.Range("E40").Value = '  +1.51%  '
# Row 41
$ws.Range("E41").Value = '  +0.00%  '
# Row 42
$ws.Range("E42").Value = '  +1.38%  '
# Row 43
$ws.Range("D43").Value = "'0.917"
$ws.Range("E43").Value = '  -0.37%  '
# Row 44
$ws.Range("D44").Value = "'36.62"
$ws.Range("E44").Value = '  +47.00%  '
# Row 45
$ws.Range("D45").Value = "'1.75"
$ws.Range("E45").Value = '  +3.23%  '
# Row 46
$ws.Range("D46").Value = "'23.67"
$ws.Range("E46").Value = '  -0.23%  '
# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = "'2.27"
$ws.Range("E47").Value = '  +6.99%  '
# Row 48
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D48").Value = "'5.66"
$ws.Range("E48").Value = '  +2.76%  '
# Row 49
$ws.Range("D49").Value = "'0.0413"
$ws.Range("E49").Value = '  -1.61%  '
# Row 50
$ws.Range("D50").Value = "'3.53"
$ws.Range("E50").Value = '  -0.90%  '
# Row 51
$ws.Range("D51").Value = "'53.52"
$ws.Range("E51").Value = '  +0.57%  '

# Restore default (General) styling on every price cell touched above so
# the quote-prefix marker introduced by the text assignment doesn't persist.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
